$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A20").Value = "BFO"
$ws.Range("C20").Value = "entity [BFO:0000001]"
$ws.Range("D20").Value = "material entity [BFO:0000040]"
$ws.Range("E20").Value = "all"
